# "Changing to caps on Nodes, etc."
# Uppercase the entity-id labels (Person1/Study1/Treat1) and lowercase the
# camelCase property/relation names (firstName/enrolledIn/treatmentArm)
# across both tables on the Neo4jModel sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string entries must land in this order: PERSON1, STUDY1,
# TREAT1, enrolledin, treatmentarm, firstname -- write cells accordingly.

# --- Table 1 (Nodes and Relations): A3:C5 -----------------------------
$ws.Range("A3").Value = "PERSON1"
$ws.Range("A4").Value = "PERSON1"
$ws.Range("E3").Value = "PERSON1"
$ws.Range("E4").Value = "PERSON1"

$ws.Range("C3").Value = "STUDY1"
$ws.Range("A5").Value = "STUDY1"
$ws.Range("E5").Value = "STUDY1"

$ws.Range("C4").Value = "TREAT1"
$ws.Range("C5").Value = "TREAT1"
$ws.Range("E6").Value = "TREAT1"
$ws.Range("E7").Value = "TREAT1"

$ws.Range("B3").Value = "enrolledin"

$ws.Range("B5").Value = "treatmentarm"

# --- Table 2 (Node P:V Pairs): E3:G7 ----------------------------------
$ws.Range("F3").Value = "firstname"

# --- Selection moved from B9 to A7 ------------------------------------
$ws.Range("A7").Select()

# --- Sheet protection enabled ------------------------------------------
$ws.Protect()
